$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.971.22'
$ws.Range("E2").Value = '  -0.37%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.826.03'
$ws.Range("E3").Value = '  -0.36%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.21'
$ws.Range("E5").Value = '  -0.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6188'
$ws.Range("E6").Value = '  -1.52%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.005'
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07382'
$ws.Range("E8").Value = '  -1.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2902'
$ws.Range("E9").Value = '  -0.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.86'
$ws.Range("E10").Value = '  -0.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07676'
$ws.Range("E11").Value = '  -0.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.829.49'
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.952'
$ws.Range("E13").Value = '  -0.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6661'
$ws.Range("E14").Value = '  -0.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.01'
$ws.Range("E15").Value = '  -0.79%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009037'
$ws.Range("E16").Value = '  -3.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.851'
$ws.Range("E17").Value = '  -2.93%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '28.989.31'
$ws.Range("E18").Value = '  -0.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.080.29'
$ws.Range("E19").Value = '  -0.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '234.65'
$ws.Range("E20").Value = '  +5.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.49'
$ws.Range("E21").Value = '  -0.72%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.006'
$ws.Range("E22").Value = '  +0.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.117'
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.006'
$ws.Range("E24").Value = '  +0.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.95'
$ws.Range("E25").Value = '  -0.81%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1407'
$ws.Range("E26").Value = '  +0.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.463'
$ws.Range("E27").Value = '  -0.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.69'
$ws.Range("E28").Value = '  -1.17%  '
$ws.Range("E29").Value = '  -0.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.100'
$ws.Range("E30").Value = '  +0.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05532'
$ws.Range("E31").Value = '  -4.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.097'
$ws.Range("E32").Value = '  -1.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.210'
$ws.Range("E33").Value = '  +0.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.832'
$ws.Range("E34").Value = '  -0.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7339'
$ws.Range("E35").Value = '  -2.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.131'
$ws.Range("E36").Value = '  -0.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.644'
$ws.Range("E37").Value = '  +0.90%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.813'
$ws.Range("E38").Value = '  +2.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01762'
$ws.Range("E39").Value = '  -1.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.197.60'
$ws.Range("E40").Value = '  -2.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.372'
$ws.Range("E41").Value = '  -2.83%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9044'
$ws.Range("E42").Value = '  +1.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.005'
$ws.Range("E43").Value = '  +0.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.86'
$ws.Range("E44").Value = '  -0.97%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.977.19'
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000122'
$ws.Range("E46").Value = '  -1.73%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '64.29'
$ws.Range("E47").Value = '  -2.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5109'
$ws.Range("E48").Value = '  +0.14%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.082'
$ws.Range("E49").Value = '  +0.71%  '
$ws.Range("B50").Value = 'TheSandbox'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4005'
$ws.Range("E50").Value = '  -1.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05804'
$ws.Range("E51").Value = '  -0.46%  '
